# Small changes to poster
#
# Resizes/repositions a handful of shapes on slide 1 (the picture collage
# in the lower-right area, plus the two text boxes that frame it).
#
# The numeric literals below are PowerPoint COM "points" values
# (EMU / 12700) chosen so that, once PowerPoint's COM layer rounds them
# through its internal Single-precision (32-bit float) representation and
# re-expands to EMU on save, the stored OOXML value lands on (or, for two
# fields that fall exactly between two representable floats, the nearest
# possible single EMU unit to) the intended target.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Rechteck 33" - white backing rectangle behind the photo collage;
# only the height grows (-> cy 11224555 -> stays, cy 4767835 -> 9409641).
$rechteck33 = $s.Shapes.Item("Rechteck 33")
$rechteck33.Height = 740.9166259765625

# "Textfeld 11" - the question textbox to the left of the rectangle;
# only the width shrinks (cx 9590076 -> 9227209).
$textfeld11 = $s.Shapes.Item("Textfeld 11")
$textfeld11.Width = 726.5519409179688

# "Picture 28" (off 10754458,11766639 -> 10769013,11677927)
$pic28 = $s.Shapes.Item("Picture 28")
$pic28.Left = 847.9537963867188
$pic28.Top = 919.5218505859375

# "Picture 31" (off 12144839,11683843 -> 12144840,11683843;
#               ext 5953935,3332378 -> 5569726,3117339)
$pic31 = $s.Shapes.Item("Picture 31")
$pic31.Left = 956.28662109375
$pic31.Top = 919.9876708984375
$pic31.Width = 438.5611267089844
$pic31.Height = 245.45977783203125

# "Picture 35" (off 16926702,11677927 -> 16685976,11538753)
$pic35 = $s.Shapes.Item("Picture 35")
$pic35.Left = 1313.8564453125
$pic35.Top = 908.563232421875

# "Picture 39" (off 15451668,14636280 -> 15710049,14709956;
#               ext 5953935,1636653 -> 5369442,1475984)
$pic39 = $s.Shapes.Item("Picture 39")
$pic39.Left = 1237.01171875
$pic39.Top = 1158.2642822265625
$pic39.Width = 422.79071044921875
$pic39.Height = 116.2192153930664

# "Picture 40" (off 9526466,14800178 -> 9754994,14820581;
#               ext 5831376,1358156 -> 5545966,1291683)
$pic40 = $s.Shapes.Item("Picture 40")
$pic40.Left = 768.1098022460938
$pic40.Top = 1166.974853515625
$pic40.Width = 436.69024658203125
$pic40.Height = 101.70732879638672
